$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Absorbance 1_01": drop the trailing (empty) column M and mark a
# handful of wells as error/missing reads ("MISSED") coming from the
# instrument's error document handling.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Absorbance 1_01")

# Remove column M entirely (was all blank cells) -> dimension becomes A1:L16
$ws1.Range("M1:M16").Delete()

# Mark a few wells as missed reads
$ws1.Range("D7").Value = "MISSED"
$ws1.Range("E8").Value = "MISSED"
$ws1.Range("C12").Value = "MISSED"

# B11 had a numeric reading (1.7455); it becomes a blank/missing text cell.
# Use the text quote-prefix trick to force an empty *text* cell (rather than
# truly clearing it, which would remove the cell altogether), then strip the
# quote-prefix formatting it introduces so the cell stays unstyled.
$ws1.Range("B11").Value = "'"
$ws1.Range("B11").ClearFormats()

# ---------------------------------------------------------------------------
# Sheet "General information": drop the trailing blank row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("General information")
$ws2.Range("A4:E4").Delete()

# ---------------------------------------------------------------------------
# Sheet "Session information": drop the trailing blank row.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Session information")
$ws3.Range("A5:E5").Delete()

# ---------------------------------------------------------------------------
# Sheet "Instrument information": drop the trailing blank row.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Instrument information")
$ws4.Range("A5:E5").Delete()

# ---------------------------------------------------------------------------
# Sheet "Layout definitions": shrink down to just columns A:B, rows 1-3.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Layout definitions")
$ws5.Range("C1:E4").Delete()
$ws5.Range("A4:B4").Delete()
